$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the combined "Ministry Course Code and Level" column into two
# columns: "Ministry Course Code" and "Ministry Course Level". Insert a
# new column after the existing course-code column (G) so the level value
# has its own column, shifting the trailing columns one place to the right.
$ws.Range("H1").EntireColumn.Insert()

# Header row
$ws.Range("G1").Value = "Ministry Course Code"
$ws.Range("H1").Value = "Ministry Course Level"

# Data rows: split "ENST 12" into code "ENST" (text) and level 12 (number)
$ws.Range("G2").Value = "ENST"
$ws.Range("H2").Value = 12
$ws.Range("G3").Value = "ENST"
$ws.Range("H3").Value = 12
$ws.Range("G4").Value = "ENST"
$ws.Range("H4").Value = 12
